# TB_Planning_Alt-Thibaud.xlsx -- "Analyse - Maquettage, choix du framework"
#
# Two new task rows are inserted into the "Planning & Journal" sheet right
# after the "Planification / Macro-planning" row (row 8):
#   - row 9  : Analyse / Maquettage, choix du framework / Terminé / 3h / 3.5h
#   - row 10 : Analyse / Mise en place de Tailwind CSS  / En cours / 1h / -
# Row 8's status also flips from "En cours" to "Terminé" now that the
# planning task is done. Everything below (the placeholder Analyse /
# Modélisation / Conception / Réalisation / Tests rows and the TOTAL row)
# shifts down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert two blank rows before row 9 (shifts old rows 9-15 to 11-17)
# ---------------------------------------------------------------------
$ws.Rows("9:10").Insert(-4121)   # xlShiftDown

# Re-apply the formatting of the row directly above (row 8) to columns A
# and F of the two new rows, so they pick up the usual border (col A) and
# the percentage "Dérive" style (col F) instead of the engine's default
# guess.
$ws.Range("A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F8").Copy()
$ws.Range("F9:F10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Row 8 : Macro-planning is now finished
# ---------------------------------------------------------------------
$ws.Range("C8").Value = "Terminé"

# ---------------------------------------------------------------------
# 3) Row 9 : new task - Maquettage, choix du framework
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Analyse"
$ws.Range("B9").Value = "Maquettage, choix du framework"
$ws.Range("C9").Value = "Terminé"
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 3.5
$ws.Range("F9").Formula = "=E9/D9"

# ---------------------------------------------------------------------
# 4) Row 10 : new task - Mise en place de Tailwind CSS
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Analyse"
$ws.Range("B10").Value = "Mise en place de Tailwind CSS"
$ws.Range("C10").Value = "En cours"
$ws.Range("D10").Value = 1
$ws.Range("F10").Formula = "=E10/D10"

# ---------------------------------------------------------------------
# 5) Fix up the AutoFilter range and the _FilterDatabase name so they
#    cover the data through the (now) row 15 instead of row 13.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("C1:F15").AutoFilter()

$filterName = $wb.Names.Item("Planning & Journal!_FilterDatabase")
$filterName.RefersTo = "='Planning & Journal'!`$C`$1:`$F`$15"

# ---------------------------------------------------------------------
# 6) Nudge the conditional-formatting ranges on column C down so they
#    keep tracking the same relative rows after the insert.
# ---------------------------------------------------------------------
$condFormats = $ws.Range("C1:C1048576").FormatConditions
$condFormats.Item(1).ModifyAppliesToRange($ws.Range("C1:C10"))
$condFormats.Item(2).ModifyAppliesToRange($ws.Range("C6:C10"))
$condFormats.Item(3).ModifyAppliesToRange($ws.Range("C1:C15"))
$condFormats.Item(4).ModifyAppliesToRange($ws.Range("C1:C1048576"))

# ---------------------------------------------------------------------
# 7) Keep the active selection sane (mirrors the saved selection in the
#    authored workbook, now pointing a bit further down the sheet).
# ---------------------------------------------------------------------
$ws.Range("C21").Select()

Write-Host "Applied: Analyse - Maquettage, choix du framework"
